$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before D (shifts old D:K data to E:L)
$ws.Range("D1").EntireColumn.Insert()

# Copy number/date formatting from column E into the newly inserted column D.
# (Rows 36/78 are fully empty and 37/79 only hold a label in column B in the
# source sheet, so they are deliberately excluded to avoid fabricating cells
# that were never there.)
$ws.Range("E7:E35").Copy()
$ws.Range("D7:D35").PasteSpecial(-4122)
$ws.Range("E38:E77").Copy()
$ws.Range("D38:D77").PasteSpecial(-4122)
$ws.Range("E80:E102").Copy()
$ws.Range("D80:D102").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# Populate new column D with the new (most-recent) period figures
$ws.Range("D7").Value = 43465
$ws.Range("D8").Value = 613400
$ws.Range("D9").Value = 174700
$ws.Range("D10").Value = 438700
$ws.Range("D12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 2400
$ws.Range("D15").Value = 197900
$ws.Range("D17").Value = 464000
$ws.Range("D18").Value = 149400
$ws.Range("D20").Value = 2100
$ws.Range("D21").Value = 349400
$ws.Range("D22").Value = 0
$ws.Range("D23").Value = 151500
$ws.Range("D24").Value = 48900
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 102600
$ws.Range("D27").Value = 102600
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = -2100
$ws.Range("D33").Value = 102600
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = 102600
$ws.Range("D38").Value = 43465
$ws.Range("D41").Value = 51000
$ws.Range("D42").Value = 32700
$ws.Range("D43").Value = 104400
$ws.Range("D44").Value = "NA"
$ws.Range("D45").Value = 14300
$ws.Range("D46").Value = 202500
$ws.Range("D47").Value = 8700
$ws.Range("D48").Value = 1312800
$ws.Range("D49").Value = 102600
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 50000
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 1676600
$ws.Range("D57").Value = 123900
$ws.Range("D58").Value = 0
$ws.Range("D59").Value = 45500
$ws.Range("D60").Value = 169400
$ws.Range("D61").Value = 399400
$ws.Range("D62").Value = 78000
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 646800
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = -298600
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 1029800
$ws.Range("D77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("D81").Value = 102600
$ws.Range("D83").Value = 197900
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 284500
$ws.Range("D91").Value = -400300
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -382600
$ws.Range("D96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = 128400
$ws.Range("D101").Value = -2700
$ws.Range("D102").Value = 27600

# Two rows (Long Term Debt / Other Liabilities) also had their newly-shifted
# second-column (E) figure corrected, not just shifted from the old D
$ws.Range("E61").Value = 256500
$ws.Range("E62").Value = 79800
